# Update the "Salida" column (B) of the "Solucion" sheet with the new
# ordering produced by the beam search / simulated annealing run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Solucion")

$newSalidas = @(
    "S012",
    "S035",
    "S031",
    "S039",
    "S026",
    "S027",
    "S037",
    "S022",
    "S018",
    "S019",
    "S023",
    "S029",
    "S040",
    "S016",
    "S025",
    "S028",
    "S010",
    "S009",
    "S017",
    "S004",
    "S032",
    "S034",
    "S003",
    "S007",
    "S030",
    "S014",
    "S008",
    "S036",
    "S013",
    "S006",
    "S011",
    "S021",
    "S033",
    "S024",
    "S001",
    "S015",
    "S038",
    "S002",
    "S005",
    "S020"
)

for ($i = 0; $i -lt $newSalidas.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newSalidas[$i]
}
